$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update existing GDP per Capita values for years 1950-2008 (rows 2-60).
# These cells store numeric-looking values as text, so we temporarily
# mark the range as Text, assign the values, then clear the formatting
# back off so the cells keep the default (unstyled) appearance.
$existingRange = $ws.Range("E2:E60")
$existingRange.NumberFormat = "@"

$ws.Range("E2").Value = "697"
$ws.Range("E3").Value = "706"
$ws.Range("E4").Value = "714"
$ws.Range("E5").Value = "722"
$ws.Range("E6").Value = "730"
$ws.Range("E7").Value = "738"
$ws.Range("E8").Value = "744"
$ws.Range("E9").Value = "751"
$ws.Range("E10").Value = "759"
$ws.Range("E11").Value = "765"
$ws.Range("E12").Value = "771"
$ws.Range("E13").Value = "779"
$ws.Range("E14").Value = "787"
$ws.Range("E15").Value = "794"
$ws.Range("E16").Value = "802"
$ws.Range("E17").Value = "810"
$ws.Range("E18").Value = "818"
$ws.Range("E19").Value = "826"
$ws.Range("E20").Value = "834"
$ws.Range("E21").Value = "842"
$ws.Range("E22").Value = "850"
$ws.Range("E23").Value = "875"
$ws.Range("E24").Value = "886"
$ws.Range("E25").Value = "944"
$ws.Range("E26").Value = "958"
$ws.Range("E27").Value = "999"
$ws.Range("E28").Value = "1023"
$ws.Range("E29").Value = "993"
$ws.Range("E30").Value = "991"
$ws.Range("E31").Value = "966"
$ws.Range("E32").Value = "1054"
$ws.Range("E33").Value = "1199"
$ws.Range("E34").Value = "1256"
$ws.Range("E35").Value = "1283"
$ws.Range("E36").Value = "1344"
$ws.Range("E37").Value = "1382"
$ws.Range("E38").Value = "1409"
$ws.Range("E39").Value = "1350"
$ws.Range("E40").Value = "1286"
$ws.Range("E41").Value = "1428"
$ws.Range("E42").Value = "1481"
$ws.Range("E43").Value = "1512.47987022857"
$ws.Range("E44").Value = "1563.90573665207"
$ws.Range("E45").Value = "1622.55451608283"
$ws.Range("E46").Value = "1718.11585629849"
$ws.Range("E47").Value = "1801.54401555344"
$ws.Range("E48").Value = "1899.89586730929"
$ws.Range("E49").Value = "2002.2359892027"
$ws.Range("E50").Value = "2050.88381062929"
$ws.Range("E51").Value = "2167.77962339752"
$ws.Range("E52").Value = "2262.23564902541"
$ws.Range("E53").Value = "2365.10333635458"
$ws.Range("E54").Value = "2486.98919138709"
$ws.Range("E55").Value = "2625.34821863314"
$ws.Range("E56").Value = "2784.68009144169"
$ws.Range("E57").Value = "2948.69961823626"
$ws.Range("E58").Value = "2896.13547876231"
$ws.Range("E59").Value = "3389.36406691748"
$ws.Range("E60").Value = "3614.47682688708"

$existingRange.ClearFormats()

# Append new rows for years 2009-2016 (rows 61-68).
$newRange = $ws.Range("E61:E68")
$newRange.NumberFormat = "@"

$ws.Range("A61").Value = 418
$ws.Range("B61").Value = "Laos"
$ws.Range("C61").Value = "GDP per Capita"
$ws.Range("D61").Value = 2009
$ws.Range("E61").Value = "3843.82974993866"

$ws.Range("A62").Value = 418
$ws.Range("B62").Value = "Laos"
$ws.Range("C62").Value = "GDP per Capita"
$ws.Range("D62").Value = 2010
$ws.Range("E62").Value = "4114.0102760815"

$ws.Range("A63").Value = 418
$ws.Range("B63").Value = "Laos"
$ws.Range("C63").Value = "GDP per Capita"
$ws.Range("D63").Value = 2011
$ws.Range("E63").Value = "4401"

$ws.Range("A64").Value = 418
$ws.Range("B64").Value = "Laos"
$ws.Range("C64").Value = "GDP per Capita"
$ws.Range("D64").Value = 2012
$ws.Range("E64").Value = "4670"

$ws.Range("A65").Value = 418
$ws.Range("B65").Value = "Laos"
$ws.Range("C65").Value = "GDP per Capita"
$ws.Range("D65").Value = 2013
$ws.Range("E65").Value = "4961"

$ws.Range("A66").Value = 418
$ws.Range("B66").Value = "Laos"
$ws.Range("C66").Value = "GDP per Capita"
$ws.Range("D66").Value = 2014
$ws.Range("E66").Value = "5251"

$ws.Range("A67").Value = 418
$ws.Range("B67").Value = "Laos"
$ws.Range("C67").Value = "GDP per Capita"
$ws.Range("D67").Value = 2015
$ws.Range("E67").Value = "5559"

$ws.Range("A68").Value = 418
$ws.Range("B68").Value = "Laos"
$ws.Range("C68").Value = "GDP per Capita"
$ws.Range("D68").Value = 2016
$ws.Range("E68").Value = "5859"

$newRange.ClearFormats()
